# Apply row-content updates per commit diff (rows 3-16 data realigned)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 94692744.0
$ws.Range("B3").Value = 78596.0
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 6462.0
$ws.Range("F3").Value = "Stuplav"
$ws.Range("G3").Value = "Nephroma bellum"
$ws.Range("H3").Value = "(Spreng.) Tuck."
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("P3").Value = "Staberget, Mpd"
$ws.Range("Q3").Value = 529542.1341322724
$ws.Range("R3").Value = 6937987.647075845

# Row 4
$ws.Range("A4").Value = 94692893.0
$ws.Range("B4").Value = 78596.0
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 6462.0
$ws.Range("F4").Value = "Stuplav"
$ws.Range("G4").Value = "Nephroma bellum"
$ws.Range("H4").Value = "(Spreng.) Tuck."
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""
$ws.Range("P4").Value = "Staberget, Mpd"
$ws.Range("Q4").Value = 529532.4741482906
$ws.Range("R4").Value = 6937931.774883721

# Row 5
$ws.Range("A5").Value = 94693842.0
$ws.Range("B5").Value = 78569.0
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458.0
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("P5").Value = "Sta, Mpd"
$ws.Range("Q5").Value = 529787.6677168193
$ws.Range("R5").Value = 6937997.680336708

# Row 6
$ws.Range("A6").Value = 94692740.0
$ws.Range("B6").Value = 78569.0
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6458.0
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("P6").Value = "Staberget, Mpd"
$ws.Range("Q6").Value = 529542.1341322724
$ws.Range("R6").Value = 6937987.647075845
$ws.Range("AC6").Value = ""

# Row 7
$ws.Range("A7").Value = 94692732.0
$ws.Range("B7").Value = 96354.0
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 221952.0
$ws.Range("F7").Value = "Spindelblomster"
$ws.Range("G7").Value = "Neottia cordata"
$ws.Range("H7").Value = "(L.) Rich."
$ws.Range("I7").Value = "15"
$ws.Range("J7").Value = "stjälkar/strån/skott"
$ws.Range("K7").Value = "överblommad"
$ws.Range("P7").Value = "Staberget, Mpd"
$ws.Range("Q7").Value = 529533.64774426
$ws.Range("R7").Value = 6938007.39613072

# Row 8
$ws.Range("A8").Value = 94692848.0
$ws.Range("B8").Value = 78596.0
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 6462.0
$ws.Range("F8").Value = "Stuplav"
$ws.Range("G8").Value = "Nephroma bellum"
$ws.Range("H8").Value = "(Spreng.) Tuck."
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("P8").Value = "Staberget, Mpd"
$ws.Range("Q8").Value = 529468.4812002254
$ws.Range("R8").Value = 6937963.017895646

# Row 9
$ws.Range("A9").Value = 94692910.0
$ws.Range("B9").Value = 103250.0
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 221725.0
$ws.Range("F9").Value = "Ögonpyrola"
$ws.Range("G9").Value = "Moneses uniflora"
$ws.Range("H9").Value = "(L.) A. Gray"
$ws.Range("I9").Value = "5"
$ws.Range("J9").Value = "stjälkar/strån/skott"
$ws.Range("K9").Value = "blomning"
$ws.Range("P9").Value = "Staberget, Mpd"
$ws.Range("Q9").Value = 529552.7958173202
$ws.Range("R9").Value = 6937931.033756874

# Row 10
$ws.Range("A10").Value = 94693866.0
$ws.Range("B10").Value = 78570.0
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 2081.0
$ws.Range("F10").Value = "Skrovellav"
$ws.Range("G10").Value = "Lobaria scrobiculata"
$ws.Range("H10").Value = "(Scop.) DC."
$ws.Range("I10").Value = ""
$ws.Range("J10").Value = ""
$ws.Range("K10").Value = ""
$ws.Range("P10").Value = "Sta, Mpd"
$ws.Range("Q10").Value = 529787.6677168193
$ws.Range("R10").Value = 6937997.680336708

# Row 11
$ws.Range("A11").Value = 94692956.0
$ws.Range("B11").Value = 78569.0
$ws.Range("D11").Value = "NT"
$ws.Range("E11").Value = 6458.0
$ws.Range("F11").Value = "Lunglav"
$ws.Range("G11").Value = "Lobaria pulmonaria"
$ws.Range("H11").Value = "(L.) Hoffm."
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("K11").Value = ""
$ws.Range("P11").Value = "Staberget, Mpd"
$ws.Range("Q11").Value = 529576.9547078039
$ws.Range("R11").Value = 6937914.19027761

# Row 12
$ws.Range("A12").Value = 94693969.0
$ws.Range("B12").Value = 77506.0
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 6425.0
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""
$ws.Range("K12").Value = ""
$ws.Range("P12").Value = "Staberget, Mpd"
$ws.Range("Q12").Value = 529829.7281196868
$ws.Range("R12").Value = 6938403.286305362

# Row 13
$ws.Range("A13").Value = 94692634.0
$ws.Range("B13").Value = 96251.0
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 219790.0
$ws.Range("F13").Value = "Fläcknycklar"
$ws.Range("G13").Value = "Dactylorhiza maculata"
$ws.Range("H13").Value = "(L.) Soó"
$ws.Range("I13").Value = ""
$ws.Range("J13").Value = ""
$ws.Range("K13").Value = ""
$ws.Range("P13").Value = "Staberget, Mpd"
$ws.Range("Q13").Value = 529547.8747119794
$ws.Range("R13").Value = 6938068.840074595

# Row 14
$ws.Range("A14").Value = 94693885.0
$ws.Range("B14").Value = 78569.0
$ws.Range("D14").Value = "NT"
$ws.Range("E14").Value = 6458.0
$ws.Range("F14").Value = "Lunglav"
$ws.Range("G14").Value = "Lobaria pulmonaria"
$ws.Range("H14").Value = "(L.) Hoffm."
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""
$ws.Range("K14").Value = ""
$ws.Range("P14").Value = "stabeget, Mpd"
$ws.Range("Q14").Value = 529791.909600144
$ws.Range("R14").Value = 6937988.036609955
$ws.Range("AC14").Value = "På gammal grov björk med uppsprucken bark"

# Row 15
$ws.Range("A15").Value = 94693081.0
$ws.Range("B15").Value = 96354.0
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 221952.0
$ws.Range("F15").Value = "Spindelblomster"
$ws.Range("G15").Value = "Neottia cordata"
$ws.Range("H15").Value = "(L.) Rich."
$ws.Range("I15").Value = "5"
$ws.Range("J15").Value = "stjälkar/strån/skott"
$ws.Range("K15").Value = "överblommad"
$ws.Range("P15").Value = "Staberget, Mpd"
$ws.Range("Q15").Value = 529632.0419017738
$ws.Range("R15").Value = 6937898.084041415

# Row 16
$ws.Range("A16").Value = 94693905.0
$ws.Range("B16").Value = 89392.0
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 1202.0
$ws.Range("F16").Value = "Ullticka"
$ws.Range("G16").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H16").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("I16").Value = ""
$ws.Range("J16").Value = ""
$ws.Range("K16").Value = ""
$ws.Range("P16").Value = "Staberget, Mpd"
$ws.Range("Q16").Value = 529799.3124251956
$ws.Range("R16").Value = 6938191.876381581
